# Update the OldCode/NewCode species-replace lookup table.
# New rows are inserted (by resorting) so the OldCode column stays
# alphabetically ordered, and several brand-new species pairs are added.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("AGROP", "AGTR"),
    @("AGRT",  "AGTR"),
    @("ANAL",  "ANAL4"),
    @("ANSE",  "ANSE4"),
    @("ARA",   "ARABI2"),
    @("ASCH",  "ASCH2"),
    @("BAMA",  "BAMA4"),
    @("BRCA",  "BRCA5"),
    @("HAFL",  "HAFL2"),
    @("HEHO",  "HEHO5"),
    @("HOBR",  "HOBR2"),
    @("LOLE",  "LOLE2"),
    @("LYDR",  "LYDR2"),
    @("MAGL",  "MAGL2")
)

$rowCount = $data.Count
$arr = New-Object 'object[,]' $rowCount, 2
for ($i = 0; $i -lt $rowCount; $i++) {
    $arr[$i, 0] = $data[$i][0]
    $arr[$i, 1] = $data[$i][1]
}

$startRow = 2
$endRow = $startRow + $rowCount - 1
$targetRange = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 2))
$targetRange.Value = $arr

$ws.Range("B7").Select()
